$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend the verification note in B14 with the extra sentence, and make it wrap like B8
$ws.Range("B14").Value = "This code is not explictly tested. It is much code, but quite trivial, and we have investigated the results throroughly when trying to understand what is happening in the model, so we deem it to be safe. The code for investigating complex II in reverse across the hypoxia range is very similar to that of GenFig1-2ABData.m."
$ws.Range("B14").WrapText = $true
$ws.Rows.Item(14).AutoFit()

# Rename function references in column A
$ws.Range("A11").Value = "genFig1_2ABData"
$ws.Range("A12").Value = "genFig3Data"

# Update the view state: scrolled down, selection on B41
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("B41").Select()
